$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("O","P","Q","R","S","T","U","V","W","X")

# --- Row 1 headers: O1:X1 = <sWord-11> .. <sWord-20> ---
$headers = @("<sWord-11>","<sWord-12>","<sWord-13>","<sWord-14>","<sWord-15>","<sWord-16>","<sWord-17>","<sWord-18>","<sWord-19>","<sWord-20>")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# --- Data rows 2-17, columns O-X ---
$data = @{
    2  = @(185,82,42,22,13,8,4,0,0,0)
    3  = @(72,26,5,1,0,0,0,0,0,0)
    4  = @(99,40,14,8,4,2,0,0,0,0)
    5  = @(189,87,60,39,16,7,1,0,0,0)
    6  = @(82,43,23,15,8,3,0,0,0,0)
    7  = @(78,47,38,28,20,12,6,3,1,0)
    8  = @(26,13,4,0,0,0,0,0,0,0)
    9  = @(16,8,3,0,0,0,0,0,0,0)
    10 = @(152,73,45,29,13,4,0,0,0,0)
    11 = @(100,45,15,5,2,0,0,0,0,0)
    12 = @(174,71,33,18,5,1,0,0,0,0)
    13 = @(266,159,90,55,26,14,5,1,0,0)
    14 = @(132,59,26,9,2,0,0,0,0,0)
    15 = @(193,105,59,23,9,5,0,0,0,0)
    16 = @(208,98,38,15,7,4,1,0,0,0)
    17 = @(106,61,33,11,4,1,0,0,0,0)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

# --- Row 19 totals: extend the shared SUM formula from column D through column X ---
for ($i = 0; $i -lt $cols.Length; $i++) {
    $c = $cols[$i]
    $ws.Range($c + "19").Formula = "=SUM(" + $c + "2:" + $c + "18)"
}

# --- Row 20 percentages: extend the shared ratio formula from column N through column X ---
for ($i = 0; $i -lt $cols.Length; $i++) {
    $c = $cols[$i]
    $ws.Range($c + "20").Formula = "=" + $c + "19/`$D`$19"
}

# --- Remove the old stray P19 = SUM(E19:L19) total (superseded by the new P19 column total above) ---
# (already overwritten by the row 19 loop above)

# --- Column widths to match the newly populated columns ---
$ws.Range("F1:J1").EntireColumn.ColumnWidth = 11.166666666666666
$ws.Range("K1").EntireColumn.ColumnWidth = 10.166666666666666
$ws.Range("L1:N1").EntireColumn.ColumnWidth = 11.166666666666666
$ws.Range("O1:S1").EntireColumn.ColumnWidth = 9.592447916666666
$ws.Range("T1:X1").EntireColumn.ColumnWidth = 10.592447916666666

# --- Freeze panes: freeze columns A:B and row 1 ---
$ws.Range("C2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Final selection in the scrollable (bottom-right) pane ---
$ws.Range("W22").Select()
